# Insert a new weekly data row at the top of the data block (row 2),
# pushing the existing rows (2-5) down to (3-6), without disturbing
# the existing cell formatting of the shifted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 5 -> 6, 4 -> 5, 3 -> 4, 2 -> 3 (bottom-up so we don't overwrite
# data before it has been copied). Columns A..T = 1..20.
# Only column D (4) carries an explicit number format (the date style), so
# only copy NumberFormat there to avoid stamping a redundant "General" style
# onto every other cell.
for ($r = 5; $r -ge 2; $r--) {
    for ($c = 1; $c -le 20; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r + 1, $c)
        $dst.Value = $src.Value2
        if ($c -eq 4) {
            $dst.NumberFormat = $src.NumberFormat
        }
    }
}

# Populate the new row 2 with the new weekly entry.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44761
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100108
$ws.Range("H2").Value = "Tropicales y subtropicales"
$ws.Range("I2").Value = 100108007
$ws.Range("J2").Value = "Coco"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("Q2").Value = "$/malla 20 unidades"
$ws.Range("R2").Value = "Perú"
$ws.Range("S2").Value = 1025
$ws.Range("T2").Value = 20
